$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row block 1: "Distance to City Center" (rows 2-3, 1-indexed) ---
# Update numeric values in row 2 (df, F, p)
$t.Cell(2, 3).Range.Text = "40"
$t.Cell(2, 5).Range.Text = "7.800"
$t.Cell(2, 6).Range.Text = "0.883"

# Update p value in row 3
$t.Cell(3, 6).Range.Text = "0.830"

# Vertically merge column 1 cells for rows 2-3 (removes duplicate
# "Distance to City Center" label from row 3, adds vMerge markup)
$t.Cell(2, 1).Merge($t.Cell(3, 1))

# --- Row block 2: "Urbanization Score" (rows 4-5, 1-indexed) ---
# Update numeric values in row 4 (df, F, p)
$t.Cell(4, 3).Range.Text = "40"
$t.Cell(4, 5).Range.Text = "6.048"
$t.Cell(4, 6).Range.Text = "0.944"

# Update p value in row 5
$t.Cell(5, 6).Range.Text = "0.944"

# Vertically merge column 1 cells for rows 4-5 (removes duplicate
# "Urbanization Score" label from row 5, adds vMerge markup)
$t.Cell(4, 1).Merge($t.Cell(5, 1))
